# Fixed some bugs in slotsgamecore7 - corrects the symbol/reel weight
# values for several reel-strip rows that were shuffled/mismatched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values for A3:F21 (symbol, reel1..reel5), row by row.
$rows = @{
    3  = @(501,  9, 52, 30, 75, 45)
    4  = @(1201, 2, 10, 10, 10, 10)
    5  = @(1202, 2, 10, 10, 10, 10)
    6  = @(901, 16, 15, 45, 60, 60)
    7  = @(401,  9, 48, 67, 75, 45)
    8  = @(301,  6, 45, 30, 60, 45)
    9  = @(701,  3, 90, 45, 97, 15)
    10 = @(801,  3, 67, 65, 52, 45)
    11 = @(902,  1,  0,  0,  0,  0)
    12 = @(1001,18, 30, 75, 60, 72)
    13 = @(1203, 3, 15, 15, 15, 15)
    14 = @(601,  9, 60, 67, 60, 42)
    15 = @(201,  9, 30, 15, 45, 30)
    16 = @(802,  0,  4,  5,  4,  0)
    17 = @(1,    0,  2,  2,  2,  2)
    18 = @(3,    0,  3,  3,  3,  3)
    19 = @(502,  0,  4,  0,  0,  0)
    20 = @(2,    0,  2,  2,  2,  2)
    21 = @(1101, 0, 15, 30, 30,  0)
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $col = [char]([int][char]'A' + $i)
        $ws.Range("$col$r").Value = $vals[$i]
    }
}
